$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width -> 21 (stored OOXML "width"). The runtime's ColumnWidth
# setter adds a fixed 5/6 (0.8333...) padding when persisting to the sheet's
# <col width=".."> attribute, so we set it 5/6 lower than the desired 21 to
# land on an on-disk width of exactly 21.
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668

$ws.Range("A2").Value = "Toyota_Yaris_2020_"
$ws.Range("A3").Value = "Mazda_MX-30_2020_"
$ws.Range("A4").Value = "Honda_Jazz_2020_"
$ws.Range("A5").Value = "Landrover Defender "
$ws.Range("A6").Value = "SEAT_Leon_2020_"
$ws.Range("A7").Value = "Kia_Sorento_2020_"
$ws.Range("A8").Value = "Honda e "
$ws.Range("A9").Value = "Hyundai_i10_2020_"
$ws.Range("A10").Value = "Isuzu_D-Max_2020_"
$ws.Range("A11").Value = "Audi A3 "
